# Applies the "Updated symbol list on Wed Jan 18 12:22:03 UTC 2023 with GitHub Actions"
# refresh to the crypto tracker sheet: per-coin price/volume/hour updates,
# a few coins re-ranking (names+links shifting rows), and the Hora column
# bumping from 11 to 12 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = 'D2'; Val = '301.21' },
    @{ Cell = 'E2'; Val = '-0.31%' },
    @{ Cell = 'G2'; Val = '12' },
    @{ Cell = 'D3'; Val = '32.29' },
    @{ Cell = 'E3'; Val = '1.70%' },
    @{ Cell = 'G3'; Val = '12' },
    @{ Cell = 'D4'; Val = '4.993' },
    @{ Cell = 'E4'; Val = '-2.22%' },
    @{ Cell = 'G4'; Val = '12' },
    @{ Cell = 'D5'; Val = '0.07635' },
    @{ Cell = 'E5'; Val = '-2.66%' },
    @{ Cell = 'G5'; Val = '12' },
    @{ Cell = 'D6'; Val = '1.964' },
    @{ Cell = 'E6'; Val = '-12.58%' },
    @{ Cell = 'G6'; Val = '12' },
    @{ Cell = 'D7'; Val = '7.834' },
    @{ Cell = 'E7'; Val = '0.38%' },
    @{ Cell = 'G7'; Val = '12' },
    @{ Cell = 'B8'; Val = 'MXToken' },
    @{ Cell = 'C8'; Val = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Cell = 'D8'; Val = '0.9191' },
    @{ Cell = 'E8'; Val = '-0.09%' },
    @{ Cell = 'G8'; Val = '12' },
    @{ Cell = 'B9'; Val = 'WazirX' },
    @{ Cell = 'C9'; Val = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' },
    @{ Cell = 'D9'; Val = '0.1755' },
    @{ Cell = 'E9'; Val = '-0.54%' },
    @{ Cell = 'G9'; Val = '12' },
    @{ Cell = 'B10'; Val = 'LiechtensteinCryptoassetsExchange' },
    @{ Cell = 'C10'; Val = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' },
    @{ Cell = 'D10'; Val = '0.07856' },
    @{ Cell = 'E10'; Val = '4.50%' },
    @{ Cell = 'G10'; Val = '12' },
    @{ Cell = 'B11'; Val = 'MandalaExchangeToken' },
    @{ Cell = 'C11'; Val = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' },
    @{ Cell = 'D11'; Val = '0.08516' },
    @{ Cell = 'E11'; Val = '-4.61%' },
    @{ Cell = 'G11'; Val = '12' },
    @{ Cell = 'B12'; Val = 'BitrueCoin' },
    @{ Cell = 'C12'; Val = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' },
    @{ Cell = 'D12'; Val = '0.03220' },
    @{ Cell = 'E12'; Val = '6.05%' },
    @{ Cell = 'G12'; Val = '12' },
    @{ Cell = 'B13'; Val = 'BitMartToken' },
    @{ Cell = 'C13'; Val = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' },
    @{ Cell = 'D13'; Val = '0.09990' },
    @{ Cell = 'E13'; Val = '-0.43%' },
    @{ Cell = 'G13'; Val = '12' },
    @{ Cell = 'B14'; Val = 'BitForexToken' },
    @{ Cell = 'C14'; Val = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' },
    @{ Cell = 'D14'; Val = '0.001510' },
    @{ Cell = 'E14'; Val = '-0.79%' },
    @{ Cell = 'G14'; Val = '12' },
    @{ Cell = 'B15'; Val = 'TigerCash' },
    @{ Cell = 'C15'; Val = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' },
    @{ Cell = 'D15'; Val = '0.005734' },
    @{ Cell = 'E15'; Val = '-0.71%' },
    @{ Cell = 'G15'; Val = '12' },
    @{ Cell = 'B16'; Val = 'UpBots' },
    @{ Cell = 'C16'; Val = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt' },
    @{ Cell = 'D16'; Val = '0.007498' },
    @{ Cell = 'E16'; Val = '2,116.77%' },
    @{ Cell = 'G16'; Val = '12' },
    @{ Cell = 'B17'; Val = 'LEO' },
    @{ Cell = 'C17'; Val = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' },
    @{ Cell = 'D17'; Val = '3.464' },
    @{ Cell = 'E17'; Val = '0.04%' },
    @{ Cell = 'G17'; Val = '12' },
    @{ Cell = 'B18'; Val = 'GateToken' },
    @{ Cell = 'C18'; Val = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' },
    @{ Cell = 'D18'; Val = '3.780' },
    @{ Cell = 'E18'; Val = '-0.91%' },
    @{ Cell = 'G18'; Val = '12' },
    @{ Cell = 'E19'; Val = '-4.41%' },
    @{ Cell = 'G19'; Val = '12' },
    @{ Cell = 'D20'; Val = '0.3339' },
    @{ Cell = 'E20'; Val = '1.41%' },
    @{ Cell = 'G20'; Val = '12' },
    @{ Cell = 'E21'; Val = '0.74%' },
    @{ Cell = 'G21'; Val = '12' },
    @{ Cell = 'D22'; Val = '4.267' },
    @{ Cell = 'E22'; Val = '4.69%' },
    @{ Cell = 'G22'; Val = '12' },
    @{ Cell = 'D23'; Val = '0.1990' },
    @{ Cell = 'E23'; Val = '9.63%' },
    @{ Cell = 'G23'; Val = '12' },
    @{ Cell = 'D24'; Val = '0.04503' },
    @{ Cell = 'E24'; Val = '-2.03%' },
    @{ Cell = 'G24'; Val = '12' },
    @{ Cell = 'E25'; Val = '-2.18%' },
    @{ Cell = 'G25'; Val = '12' },
    @{ Cell = 'D26'; Val = '0.004402' },
    @{ Cell = 'E26'; Val = '-1.66%' },
    @{ Cell = 'G26'; Val = '12' },
    @{ Cell = 'D27'; Val = '0.0001250' },
    @{ Cell = 'E27'; Val = '0.23%' },
    @{ Cell = 'G27'; Val = '12' },
    @{ Cell = 'G28'; Val = '12' },
    @{ Cell = 'G29'; Val = '12' },
    @{ Cell = 'G30'; Val = '12' },
    @{ Cell = 'G31'; Val = '12' },
    @{ Cell = 'G32'; Val = '12' },
    @{ Cell = 'G33'; Val = '12' },
    @{ Cell = 'G34'; Val = '12' },
    @{ Cell = 'G35'; Val = '12' },
    @{ Cell = 'G36'; Val = '12' },
    @{ Cell = 'G37'; Val = '12' },
    @{ Cell = 'G38'; Val = '12' },
    @{ Cell = 'D39'; Val = '0.01710' },
    @{ Cell = 'E39'; Val = '-3.48%' },
    @{ Cell = 'G39'; Val = '12' },
    @{ Cell = 'D40'; Val = '0.04677' },
    @{ Cell = 'E40'; Val = '-2.41%' },
    @{ Cell = 'G40'; Val = '12' },
    @{ Cell = 'E41'; Val = '1.81%' },
    @{ Cell = 'G41'; Val = '12' },
    @{ Cell = 'D42'; Val = '0.1349' },
    @{ Cell = 'E42'; Val = '-1.02%' },
    @{ Cell = 'G42'; Val = '12' },
    @{ Cell = 'E43'; Val = '6.63%' },
    @{ Cell = 'G43'; Val = '12' },
    @{ Cell = 'D44'; Val = '0.01050' },
    @{ Cell = 'E44'; Val = '0.17%' },
    @{ Cell = 'G44'; Val = '12' },
    @{ Cell = 'D45'; Val = '0.00006261' },
    @{ Cell = 'E45'; Val = '-0.24%' },
    @{ Cell = 'G45'; Val = '12' },
    @{ Cell = 'E46'; Val = '0.16%' },
    @{ Cell = 'G46'; Val = '12' },
    @{ Cell = 'D47'; Val = '0.003000' },
    @{ Cell = 'E47'; Val = '-62.44%' },
    @{ Cell = 'G47'; Val = '12' },
    @{ Cell = 'D48'; Val = '0.8204' },
    @{ Cell = 'E48'; Val = '-26.33%' },
    @{ Cell = 'G48'; Val = '12' },
    @{ Cell = 'D49'; Val = '0.00002100' },
    @{ Cell = 'E49'; Val = '0.16%' },
    @{ Cell = 'G49'; Val = '12' },
    @{ Cell = 'D50'; Val = '0.0002000' },
    @{ Cell = 'E50'; Val = '0.16%' },
    @{ Cell = 'G50'; Val = '12' },
    @{ Cell = 'G51'; Val = '12' }
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    $cell.Value = "'" + $change.Val
    $cell.Style = "Normal"
}

Write-Output ("Applied {0} cell updates" -f $changes.Count)
